# Saldo.xlsx update:
#  - LAILA's balance (row 2) changes from 318117.54 to 350000
#  - LEVI (005206566) row is removed
#  - LUISA (004855570) row (the one with balance 12526.14) is removed
#  - KELMA (004504449) moves from the bottom of the list (balance 13.75)
#    to just above FABIOLA, with a new balance of 1013.75

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) LAILA's balance: 318117.54 -> 350000
$ws.Range("C2").Value = 350000

# 2) Remove the LEVI row (account 005206566), originally row 5.
$ws.Rows(5).Delete()

# 3) Remove the LUISA row (account 004855570), originally row 9 -
#    now row 8 since the LEVI row above it was already deleted.
$ws.Rows(8).Delete()

# 4) Remove KELMA's old row (account 004504449, balance 13.75),
#    originally row 363 - now row 361 after the two deletions above.
$ws.Rows(361).Delete()

# 5) Insert KELMA's new row just above FABIOLA. FABIOLA was originally
#    row 16, now row 14 after the two row deletions above.
$ws.Rows(14).Insert()

# Keep the leading zeros in the account number (text, not a number):
# format as text, enter the value, then drop back to the default/no
# explicit format so the new row matches its neighbours (which also
# carry no explicit cell style).
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "004504449"
$ws.Range("A14").ClearFormats()

$ws.Range("B14").Value = "KELMA"
$ws.Range("C14").Value = 1013.75
